$d = $word.ActiveDocument

# Locate the target paragraph (the one beginning "You are participating ...")
# by scanning the document's paragraphs rather than hard-coding an index.
$paras = $d.Paragraphs
$target = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $cand = $paras.Item($i)
    if ($cand.Range.Text -like "You are participating*") {
        $target = $cand
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the 'You are participating ...' paragraph"
}

$newText = "You are participating in a global campaign to observe and record the faintest stars visible as a means of measuring light pollution in a given location. By locating and observing the constellation Orion constellation in the night sky and comparing it to stellar charts, people from around the world will learn how the lights in their community contribute to light pollution. Your contributions to the online database will document the visible nighttime sky."

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="1312125B" w14:textId="3C7DDB2E" w:rsidR="004615A9" w:rsidRPr="00DB0F3B" w:rsidRDefault="00BE6DBA" w:rsidP="004615A9"><w:pPr><w:pStyle w:val="BasicParagraph"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="-72"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Optima-Regular"/><w:sz w:val="20"/></w:rPr></w:pPr><w:r><w:t>' + $newText + '</w:t></w:r></w:p>'

$target.Range.InsertXML($xml)
